# Update column C ("Förändrad" date) for rows 2-151: 45189 -> 45190 (+1 day)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 151; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45189) {
        $cell.Value2 = 45190
    }
}
